$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4.618035477178182,
    8.121859621296721,
    7.769047614292191,
    8.851964017779409,
    0.5330153452348441,
    8.707326061089759,
    4.765926304770545,
    3.412055793969955,
    6.86535470680758,
    7.697127969469818,
    0.2141351277505805,
    4.398894272988684,
    2.454402358756539,
    4.970242716006027,
    4.773664190038545,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231,
    8.524791795037231
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
